$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D price cells keep their exact text representation
# (avoid Excel auto-converting numeric-looking strings to numbers,
# which would strip meaningful trailing zeros like "587.00" -> 587).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '69.689.18'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '3.521.45'
$ws.Range('E3').Value = '  -2.33%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '587.00'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').Value = '183.35'
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('D7').Value = '3.511.33'
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  -2.81%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('D10').Value = '0.196'
$ws.Range('E10').Value = '  +5.21%  '
$ws.Range('D11').Value = '0.643'
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('D12').Value = '54.06'
$ws.Range('E12').Value = '  -3.30%  '
$ws.Range('D13').Value = '0.0000303'
$ws.Range('E13').Value = '  -2.83%  '
$ws.Range('D14').Value = '9.45'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '4.080.76'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('D16').Value = '19.30'
$ws.Range('E16').Value = '  -2.45%  '
$ws.Range('D17').Value = '69.696.42'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '3.507.97'
$ws.Range('E18').Value = '  -2.48%  '
$ws.Range('D19').Value = '12.36'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('E20').Value = '  -1.23%  '
$ws.Range('D21').Value = '543.18'
$ws.Range('E21').Value = '  +11.10%  '
$ws.Range('E22').Value = '  -3.29%  '
$ws.Range('D23').Value = '17.80'
$ws.Range('E23').Value = '  -8.28%  '
$ws.Range('D24').Value = '4.55'
$ws.Range('E24').Value = '  +4.58%  '
$ws.Range('D25').Value = '4.83'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('D26').Value = '95.68'
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = '11.13'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('B28').Value = 'ImmutableX'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D28').Value = '2.97'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').Value = '9.05'
$ws.Range('E29').Value = '  -3.48%  '
$ws.Range('D30').Value = '32.04'
$ws.Range('E30').Value = '  -0.46%  '
$ws.Range('D31').Value = '7.24'
$ws.Range('E31').Value = '  -4.85%  '
$ws.Range('D32').Value = '12.39'
$ws.Range('D33').Value = '64.28'
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('E34').Value = '  -4.02%  '
$ws.Range('D35').Value = '545.65'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').Value = '0.409'
$ws.Range('E36').Value = '  +3.47%  '
$ws.Range('E37').Value = '  +4.34%  '
$ws.Range('D38').Value = '37.94'
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('D40').Value = '0.0₃0758'
$ws.Range('E40').Value = '  -6.53%  '
$ws.Range('D41').Value = '3.367.98'
$ws.Range('E41').Value = '  +4.68%  '
$ws.Range('E42').Value = '  -2.47%  '
$ws.Range('D43').Value = '3.37'
$ws.Range('E43').Value = '  -3.56%  '
$ws.Range('D44').Value = '3.08'
$ws.Range('E44').Value = '  -7.41%  '
$ws.Range('D45').Value = '3.55'
$ws.Range('E45').Value = '  +0.57%  '
$ws.Range('D46').Value = '2.96'
$ws.Range('E46').Value = '  -1.64%  '
$ws.Range('D47').Value = '0.0438'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '9.09'
$ws.Range('E48').Value = '  -6.28%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.135'
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = '136.70'
$ws.Range('E51').Value = '  +1.55%  '
